# Add a new "2021" column (column R) to the sheet: one more year of data
# to the right of the existing "2020" column (Q). Each new R cell should
# carry the same formatting as the Q cell in its row, so we copy Q's
# format onto R first and then overwrite the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the year headers (...,2019,2020,2021,...).
# Rows 4-33 hold the per-region/per-sex data values for 2021.
$values = [ordered]@{
    3  = 2021
    4  = 1.7931687443515183
    5  = 1.0977143806517458
    6  = 2.4989281705678046
    7  = 2.3489023398681002
    8  = 1.8410239038543676
    9  = 2.8382683724659588
    10 = 1.2584206034913306
    11 = 0.79202525610136665
    12 = 1.7183687369364922
    13 = 1.7860084101151579
    14 = 1.5807090270340762
    15 = 1.9930959157478496
    16 = 1.0231016349164126
    17 = 0
    18 = 2.0091214112068791
    19 = 2.2092990108041848
    20 = 0.86496336159360854
    21 = 3.5236628052020538
    22 = 1.4678252700798498
    23 = 0.74155920237892192
    24 = 2.1792664589099311
    25 = 1.5302890103825006
    26 = 0.80351618683358383
    27 = 2.280288974802807
    28 = 2.3014726663297309
    29 = 1.7358308467556451
    30 = 2.9402079315049163
    31 = 1.2198989923634325
    32 = 1.1878318505232399
    33 = 1.2537455648750642
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    # Copy the formatting (number format, font, border, ...) from the
    # existing Q column cell onto the new R column cell for this row.
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = 0

    # Write the new 2021 value into the cell.
    $dstCell.Value = $values[$row]
}

# Match the post-edit selection captured for this sheet.
$ws.Range("S14").Select()
